$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H11").Value = 114.4
$ws.Range("I11").Value = 114.4
$ws.Range("K11").Value = 114.4
$ws.Range("M11").Value = 25.59999999999999

$ws.Range("H86").Value = 2381
$ws.Range("I86").Value = 2115.7144
$ws.Range("J86").Value = 3000
$ws.Range("K86").Value = 2115.7144
$ws.Range("L86").Value = 3000
$ws.Range("M86").Value = -992.7143999999998
$ws.Range("N86").Value = -5246

$ws.Range("H89").Value = 2381
$ws.Range("I89").Value = 2115.7144
$ws.Range("J89").Value = 3000
$ws.Range("K89").Value = 10578.572
$ws.Range("L89").Value = 15000
$ws.Range("M89").Value = -4962.572
$ws.Range("N89").Value = -26232

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H37").Value = 32900
$ws.Range("I37").Value = 0
$ws.Range("J37").Value = 32900
$ws.Range("K37").Value = 0
$ws.Range("L37").Value = 32900
$ws.Range("N37").Value = -33446
$ws.Range("M37").ClearContents()

$ws.Range("H97").Value = 821.2308
$ws.Range("I97").Value = 848.7
$ws.Range("J97").Value = 729.6667
$ws.Range("K97").Value = 848.7
$ws.Range("L97").Value = 729.6667
$ws.Range("M97").Value = -352.7
$ws.Range("N97").Value = -1721.6667

$ws.Range("H110").Value = 2739.875
$ws.Range("J110").Value = 3355
$ws.Range("L110").Value = 3355
$ws.Range("N110").Value = -7445

$ws.Range("H122").Value = 1198.1818
$ws.Range("I122").Value = 1107.3334
$ws.Range("J122").Value = 1607
$ws.Range("K122").Value = 3322.0002
$ws.Range("L122").Value = 4821
$ws.Range("M122").Value = -872.0001999999999
$ws.Range("N122").Value = -9721

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 146944.58
$ws.Range("I86").Value = 5435
$ws.Range("J86").Value = 253076.75
$ws.Range("K86").Value = 5435
$ws.Range("L86").Value = 253076.75
$ws.Range("M86").Value = -4312
$ws.Range("N86").Value = -255322.75

$ws.Range("H89").Value = 146944.58
$ws.Range("I89").Value = 5435
$ws.Range("J89").Value = 253076.75
$ws.Range("K89").Value = 27175
$ws.Range("L89").Value = 1265383.75
$ws.Range("M89").Value = -21559
$ws.Range("N89").Value = -1276615.75

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H99").Value = 1663.8572
$ws.Range("I99").Value = 1798.1818
$ws.Range("J99").Value = 1171.3334
$ws.Range("K99").Value = 1798.1818
$ws.Range("L99").Value = 1171.3334
$ws.Range("M99").Value = -300.1818000000001
$ws.Range("N99").Value = -4167.3334

$ws.Range("H122").Value = 2900.75
$ws.Range("I122").Value = 2081.1428
$ws.Range("K122").Value = 6243.428400000001
$ws.Range("M122").Value = -3793.428400000001

$ws.Range("H126").Value = 1663.8572
$ws.Range("I126").Value = 1798.1818
$ws.Range("J126").Value = 1171.3334
$ws.Range("K126").Value = 5394.5454
$ws.Range("L126").Value = 3514.0002
$ws.Range("M126").Value = -2924.5454
$ws.Range("N126").Value = -8454.0002

$ws.Range("H134").Value = 1872.7693
$ws.Range("I134").Value = 1497.1852
$ws.Range("J134").Value = 2717.8333
$ws.Range("K134").Value = 4491.5556
$ws.Range("L134").Value = 8153.499899999999
$ws.Range("M134").Value = -1956.5556
$ws.Range("N134").Value = -13223.4999

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H40").Value = 530.4545000000001
$ws.Range("I40").Value = 158.33333
$ws.Range("J40").Value = 670
$ws.Range("K40").Value = 633.33332
$ws.Range("L40").Value = 2680
$ws.Range("M40").Value = -564.33332
$ws.Range("N40").Value = -2818

$ws.Range("H117").Value = 301.4
$ws.Range("I117").Value = 262.33334
$ws.Range("J117").Value = 360
$ws.Range("K117").Value = 787.0000200000001
$ws.Range("L117").Value = 1080
$ws.Range("M117").Value = 2654.99998
$ws.Range("N117").Value = -7964

$ws.Range("H131").Value = 27779158
$ws.Range("I131").Value = 482.85715
$ws.Range("J131").Value = 34484356
$ws.Range("K131").Value = 1448.57145
$ws.Range("L131").Value = 103453068
$ws.Range("M131").Value = 3591.42855
$ws.Range("N131").Value = -103463148

$ws.Range("H132").Value = 2154.6
$ws.Range("I132").Value = 0
$ws.Range("J132").Value = 2154.6
$ws.Range("K132").Value = 0
$ws.Range("L132").Value = 19391.4
$ws.Range("N132").Value = -24451.4
$ws.Range("M132").ClearContents()

$ws.Range("H140").Value = 1826.5652
$ws.Range("I140").Value = 857.875
$ws.Range("K140").Value = 2573.625
$ws.Range("M140").Value = 2606.375

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H102").Value = 2696
$ws.Range("I102").Value = 2647.111
$ws.Range("J102").Value = 2758.8572
$ws.Range("K102").Value = 2647.111
$ws.Range("L102").Value = 2758.8572
$ws.Range("M102").Value = -1025.111
$ws.Range("N102").Value = -6002.8572

$ws.Range("H122").Value = 3003.9048
$ws.Range("I122").Value = 2382.182
$ws.Range("K122").Value = 7146.545999999999
$ws.Range("M122").Value = -4696.545999999999

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H22").Value = 60007.5
$ws.Range("I22").Value = 50000
$ws.Range("J22").Value = 70015
$ws.Range("K22").Value = 50000
$ws.Range("L22").Value = 70015
$ws.Range("M22").Value = -49707
$ws.Range("N22").Value = -70601

$ws.Range("H70").Value = 49909.09
$ws.Range("J70").Value = 16555.555
$ws.Range("L70").Value = 16555.555
$ws.Range("N70").Value = -17185.555

$ws.Range("H73").Value = 49909.09
$ws.Range("J73").Value = 16555.555
$ws.Range("L73").Value = 16555.555
$ws.Range("N73").Value = -18739.555

$ws.Range("H81").Value = 70534.3
$ws.Range("I81").Value = 126420.336
$ws.Range("J81").Value = 7662.5
$ws.Range("K81").Value = 252840.672
$ws.Range("L81").Value = 15325
$ws.Range("M81").Value = -251779.672
$ws.Range("N81").Value = -17447

$ws.Range("H84").Value = 70534.3
$ws.Range("I84").Value = 126420.336
$ws.Range("J84").Value = 7662.5
$ws.Range("K84").Value = 1264203.36
$ws.Range("L84").Value = 76625
$ws.Range("M84").Value = -1258899.36
$ws.Range("N84").Value = -87233

$ws.Range("H107").Value = 813.1429000000001
$ws.Range("I107").Value = 760.4
$ws.Range("J107").Value = 945
$ws.Range("K107").Value = 2281.2
$ws.Range("L107").Value = 2835
$ws.Range("M107").Value = -361.1999999999998
$ws.Range("N107").Value = -6675

$ws.Range("H109").Value = 36427
$ws.Range("J109").Value = 36664.832
$ws.Range("L109").Value = 36664.832
$ws.Range("N109").Value = -39438.832

$ws.Range("H126").Value = 4505.1177
$ws.Range("I126").Value = 4360.846
$ws.Range("K126").Value = 13082.538
$ws.Range("M126").Value = -10612.538

$ws.Range("H132").Value = 3371.8965
$ws.Range("J132").Value = 3426.6667
$ws.Range("L132").Value = 10280.0001
$ws.Range("N132").Value = -15340.0001

$ws.Range("H136").Value = 1917.7556
$ws.Range("I136").Value = 1696.875
$ws.Range("K136").Value = 5090.625
$ws.Range("M136").Value = -2540.625
